# Updates the Weekly Planner sheet:
#  - inserts a missing Quiz 2 / Assignment 3 entry (shifting the
#    Assignment/Quiz numbering on later rows down by one)
#  - reworks the final-project wrap-up rows into iteration 3 / iteration 4 /
#    demonstration milestones
#  - widens column D to fit the new (longer) text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Renumber the Assignment sequence (shifted by one) ---
$ws.Range("D16").Value = "Assignment 7"
$ws.Range("D21").Value = "Assignment 8"
$ws.Range("D24").Value = "Assignment 9"
$ws.Range("D27").Value = "Assignment 10"

# --- Class 3 / Class 4 block: add the missing Quiz 2 + Assignment 3 entries ---
$ws.Range("E4").Value = "Quiz 2"
$ws.Range("D5").Value = "Assignment 3"
$ws.Range("E5").Value = "Quiz 3"

# --- Renumber the Quiz sequence that follows (shifted by one) ---
$ws.Range("D6").Value = "Assignment 4"
$ws.Range("E7").Value = "Quiz 4"
$ws.Range("E9").Value = "Quiz 5"
$ws.Range("E11").Value = "Quiz 6"

$ws.Range("D12").Value = "Assignment 6"
$ws.Range("E12").Value = " "

# --- Final project milestones rework ---
$ws.Range("D37").Value = "Final Project - iteration 4"
$ws.Range("D38").Value = "Final Project - demonstration"

# --- Column D is now wider to fit the longer assignment/milestone text ---
$ws.Columns("D").ColumnWidth = 26.5546875

# --- Last active selection in the sheet when it was saved ---
$ws.Range("J17").Select()
